# Insert a new "location" column between "site" (A) and "lat" (B),
# shifting the existing lat/long columns from B:C to C:D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").EntireColumn.Insert()

$ws.Range("B1").Value = "location "
$ws.Range("B2:B4").Value = "anacapa_island"

$ws.Range("B5").Select()
